$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览): update "想去人数" (F column) values ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 530
$ws1.Range("F4").Value = 619
$ws1.Range("F6").Value = 454
$ws1.Range("F7").Value = 37
$ws1.Range("F8").Value = 2109
$ws1.Range("F9").Value = 851
$ws1.Range("F10").Value = 818
$ws1.Range("F12").Value = 70
$ws1.Range("F13").Value = 420
$ws1.Range("F14").Value = 319
$ws1.Range("F16").Value = 862
$ws1.Range("F18").Value = 26
$ws1.Range("F19").Value = 1647
$ws1.Range("F20").Value = 43
$ws1.Range("F21").Value = 27
$ws1.Range("F22").Value = 24
$ws1.Range("F23").Value = 58
$ws1.Range("F25").Value = 1451
$ws1.Range("F26").Value = 10
$ws1.Range("F27").Value = 517
$ws1.Range("F28").Value = 347
$ws1.Range("F29").Value = 572
$ws1.Range("F30").Value = 407
$ws1.Range("F31").Value = 2300
$ws1.Range("F32").Value = 380
$ws1.Range("F33").Value = 83
$ws1.Range("F34").Value = 165
$ws1.Range("F35").Value = 591
$ws1.Range("F36").Value = 461
$ws1.Range("F37").Value = 178
$ws1.Range("F38").Value = 908
$ws1.Range("F41").Value = 394
$ws1.Range("F42").Value = 359

# --- Sheet 2 (演出): remove the ALEXANDROS row (row 3), shifting rows 4-25 up ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(3).Delete()

# Restore sequential index values in column A (Excels row delete shifts them, but source data keeps them fixed)
for ($r = 3; $r -le 24; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

# The want-to-go count for the last (now-shifted) row increased independently
$ws2.Range("F24").Value = 429

# --- Sheet 3 (本地生活): update F column values ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 2914
$ws3.Range("F4").Value = 372
$ws3.Range("F5").Value = 240
$ws3.Range("F6").Value = 295

# --- Sheet 4 (全部类型): update F column values ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 530
$ws4.Range("F6").Value = 372
$ws4.Range("F7").Value = 240
$ws4.Range("F8").Value = 619
$ws4.Range("F11").Value = 454
$ws4.Range("F12").Value = 2109
$ws4.Range("F13").Value = 851
$ws4.Range("F14").Value = 818
$ws4.Range("F16").Value = 70
$ws4.Range("F17").Value = 420
$ws4.Range("F18").Value = 319
$ws4.Range("F20").Value = 863
$ws4.Range("F21").Value = 26
$ws4.Range("F22").Value = 295
$ws4.Range("F23").Value = 1647
$ws4.Range("F24").Value = 43
$ws4.Range("F27").Value = 58
$ws4.Range("F31").Value = 1451
$ws4.Range("F33").Value = 10
$ws4.Range("F34").Value = 517
$ws4.Range("F35").Value = 572
$ws4.Range("F36").Value = 407
$ws4.Range("F38").Value = 2300
$ws4.Range("F39").Value = 83
$ws4.Range("F40").Value = 165
$ws4.Range("F41").Value = 591
$ws4.Range("F42").Value = 461
$ws4.Range("F43").Value = 178
$ws4.Range("F44").Value = 908
$ws4.Range("F46").Value = 86
$ws4.Range("F47").Value = 429
